$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per Dr Hou's advice: expand Calr-Itga3 results to the full 4x4 grid of
# Sending cluster x Target cluster combinations (ECs, FAPs, M2, sCs),
# updating existing rows 2-9 and adding new rows 10-17.
$data = @(
    ,@("ECs", "Calr", "Itga3", "ECs", 3, 1, 81.03634266666667, 243.109028, 0.1632931649012984, 0.1632931649012984, 2, 0.6666666666666666, 5.970993, 17.912979, 0.6157237531330177, 0.6157237531330177, 483.867434808268, 4354.806913274412, 0.1005434803539962, 0.1005434803539962)
    ,@("ECs", "Calr", "Itga3", "FAPs", 3, 1, 81.03634266666667, 243.109028, 0.1632931649012984, 0.1632931649012984, 3, 1, 1.061748, 3.185244, 0.10948655666511, 0.10948655666511, 86.040174753648, 774.361572782832, 0.01787840635199116, 0.01787840635199116)
    ,@("ECs", "Calr", "Itga3", "M2", 3, 1, 81.03634266666667, 243.109028, 0.1632931649012984, 0.1632931649012984, 1, 0.3333333333333333, 0.02952066666666667, 0.088562, 0.003044146203987976, 0.003044146203987975, 2.392246859748445, 21.530221737736, 0.0004970882680714699, 0.0004970882680714699)
    ,@("ECs", "Calr", "Itga3", "sCs", 3, 1, 81.03634266666667, 243.109028, 0.1632931649012984, 0.1632931649012984, 3, 1, 2.635257666666666, 7.905773, 0.2717455439978843, 0.2717455439978843, 213.5516432909604, 1921.964789618644, 0.04437418992723955, 0.04437418992723954)
    ,@("FAPs", "Calr", "Itga3", "ECs", 3, 1, 178.365814, 535.097442, 0.3594179761796791, 0.3594179761796791, 2, 0.6666666666666666, 5.970993, 17.912979, 0.6157237531330177, 0.6157237531330177, 1065.021026833302, 9585.189241499718, 0.2213021852368256, 0.2213021852368256)
    ,@("FAPs", "Calr", "Itga3", "FAPs", 3, 1, 178.365814, 535.097442, 0.3594179761796791, 0.3594179761796791, 3, 1, 1.061748, 3.185244, 0.10948655666511, 0.10948655666511, 189.379546282872, 1704.415916545848, 0.0393514366154556, 0.0393514366154556)
    ,@("FAPs", "Calr", "Itga3", "M2", 3, 1, 178.365814, 535.097442, 0.3594179761796791, 0.3594179761796791, 1, 0.3333333333333333, 0.02952066666666667, 0.088562, 0.003044146203987976, 0.003044146203987975, 5.265477739822667, 47.389299658404, 0.001094120867832411, 0.001094120867832411)
    ,@("FAPs", "Calr", "Itga3", "sCs", 3, 1, 178.365814, 535.097442, 0.3594179761796791, 0.3594179761796791, 3, 1, 2.635257666666666, 7.905773, 0.2717455439978843, 0.2717455439978843, 470.0398788147406, 4230.358909332666, 0.09767023345956552, 0.09767023345956549)
    ,@("M2", "Calr", "Itga3", "ECs", 3, 1, 146.2303303333333, 438.6909910000001, 0.2946630198121519, 0.2946630198121519, 2, 0.6666666666666666, 5.970993, 17.912979, 0.6157237531330177, 0.6157237531330177, 873.1402788080211, 7858.26250927219, 0.1814310204682469, 0.1814310204682469)
    ,@("M2", "Calr", "Itga3", "FAPs", 3, 1, 146.2303303333333, 438.6909910000001, 0.2946630198121519, 0.2946630198121519, 3, 1, 1.061748, 3.185244, 0.10948655666511, 0.10948655666511, 155.259760770756, 1397.337846936804, 0.03226163941577561, 0.03226163941577561)
    ,@("M2", "Calr", "Itga3", "M2", 3, 1, 146.2303303333333, 438.6909910000001, 0.2946630198121519, 0.2946630198121519, 1, 0.3333333333333333, 0.02952066666666667, 0.088562, 0.003044146203987976, 0.003044146203987975, 4.316816838326889, 38.85135154494201, 0.0008969973132167959, 0.0008969973132167958)
    ,@("M2", "Calr", "Itga3", "sCs", 3, 1, 146.2303303333333, 438.6909910000001, 0.2946630198121519, 0.2946630198121519, 3, 1, 2.635257666666666, 7.905773, 0.2717455439978843, 0.2717455439978843, 385.3545991101159, 3468.191391991043, 0.08007336261491259, 0.08007336261491257)
    ,@("sCs", "Calr", "Itga3", "ECs", 3, 1, 90.630432, 271.891296, 0.1826258391068707, 0.1826258391068707, 2, 0.6666666666666666, 5.970993, 17.912979, 0.6157237531330177, 0.6157237531330177, 541.153675058976, 4870.383075530784, 0.1124470670739491, 0.1124470670739491)
    ,@("sCs", "Calr", "Itga3", "FAPs", 3, 1, 90.630432, 271.891296, 0.1826258391068707, 0.1826258391068707, 3, 1, 1.061748, 3.185244, 0.10948655666511, 0.10948655666511, 96.22667991513599, 866.040119236224, 0.01999507428188767, 0.01999507428188766)
    ,@("sCs", "Calr", "Itga3", "M2", 3, 1, 90.630432, 271.891296, 0.1826258391068707, 0.1826258391068707, 1, 0.3333333333333333, 0.02952066666666667, 0.088562, 0.003044146203987976, 0.003044146203987975, 2.675470772928, 24.079236956352, 0.0005559397548672992, 0.0005559397548672991)
    ,@("sCs", "Calr", "Itga3", "sCs", 3, 1, 90.630432, 271.891296, 0.1826258391068707, 0.1826258391068707, 3, 1, 2.635257666666666, 7.905773, 0.2717455439978843, 0.2717455439978843, 238.834540761312, 2149.510866851808, 0.04962775799616667, 0.04962775799616666)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value2 = $rowData[$c]
    }
}